$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean)
$ws.Range("B3").Value = 5274.753221354645
$ws.Range("D3").Value = 358.7562405079912
$ws.Range("E3").Value = 1026.394109535767

# Row 4 (std)
$ws.Range("B4").Value = 2239.107111578818
$ws.Range("D4").Value = 260.2557971754076
$ws.Range("E4").Value = 307.0650882933824

# Row 5 (min)
$ws.Range("B5").Value = 1572.434019178082
$ws.Range("D5").Value = 4.657665753424657
$ws.Range("E5").Value = 408.6635041095888

# Row 6 (25%)
$ws.Range("B6").Value = 3455.605529452056
$ws.Range("D6").Value = 93.151101369863
$ws.Range("E6").Value = 763.0724561643824

# Row 7 (50%)
$ws.Range("B7").Value = 4550.077228767123
$ws.Range("D7").Value = 319.8356164383562
$ws.Range("E7").Value = 1048.694015068493

# Row 8 (75%)
$ws.Range("B8").Value = 7141.529674657531
$ws.Range("D8").Value = 640
$ws.Range("E8").Value = 1323.794597260274

# Row 9 (max)
$ws.Range("B9").Value = 10047.15101095893
$ws.Range("D9").Value = 644.6575863013696
$ws.Range("E9").Value = 1676.743756164386

# Row 10 (Total)
$ws.Range("F10").Value = 7595644.638750695

# Row 11 (Residential)
$ws.Range("G11").Value = 0.7373999707823244

# Row 12 (Community)
$ws.Range("F12").Value = 516608.9863315067
$ws.Range("G12").Value = 0.06801384357766332

# Row 13 (IGA)
$ws.Range("F13").Value = 1478007.517731505
$ws.Range("G13").Value = 0.1945861856400121
